$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: replace BJ3 value with BK3
$ws.Range("BJ3").ClearContents()
$ws.Range("BK3").Value = 364.056351771112

# Row 6: remove BF6
$ws.Range("BF6").ClearContents()

# Rows 8-30: fix column A index values (shift down by 1 starting row 8)
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10
$ws.Range("A13").Value = 11
$ws.Range("A14").Value = 12
$ws.Range("A15").Value = 13
$ws.Range("A16").Value = 14
$ws.Range("A17").Value = 15
$ws.Range("A18").Value = 16
$ws.Range("A19").Value = 17
$ws.Range("A20").Value = 18
$ws.Range("A21").Value = 19
$ws.Range("A22").Value = 20
$ws.Range("A23").Value = 21
$ws.Range("A24").Value = 22
$ws.Range("A25").Value = 23
$ws.Range("A26").Value = 24
$ws.Range("A27").Value = 25
$ws.Range("A28").Value = 26
$ws.Range("A29").Value = 27
$ws.Range("A30").Value = 28

# Row 9: remove BH9, BL9
$ws.Range("BH9").ClearContents()
$ws.Range("BL9").ClearContents()

# Row 13: remove BH13, BL13
$ws.Range("BH13").ClearContents()
$ws.Range("BL13").ClearContents()

# Row 17: remove BH17, BL17
$ws.Range("BH17").ClearContents()
$ws.Range("BL17").ClearContents()

# Row 19: remove BH19, BL19
$ws.Range("BH19").ClearContents()
$ws.Range("BL19").ClearContents()

# Row 23: add BK23
$ws.Range("BK23").Value = 111.0799596430522

# Row 29: remove BH29, BL29
$ws.Range("BH29").ClearContents()
$ws.Range("BL29").ClearContents()
